# Commit: "added new test cases"
# - Sheet1!A2 email address was corrected from a gmail.com to a yahoo.com address.
# - The active selection on Sheet1 moved from I6 to D3.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the email address stored in Sheet1!A2 (shared string table entry
# "gayatri@gmail.com" -> "gayatri@yahoo.com").
$ws.Range("A2").Value = "gayatri@yahoo.com"

# Make sure we're working on the sheet the selection change applies to, then
# move the active cell / selection to D3 (was I6).
$ws.Activate()
$ws.Range("D3").Select()
